# Vendors.xlsx: normalize "Stock" column to 3 across all three sheets
# (Grocery, Pet, Bath), except the "Catnip" row on Pet which keeps its
# existing value (it isn't recorded by the same sequence/pass).
# Also move the active-sheet/selection state: Bath becomes the selected
# tab (with B6 selected), Grocery's selection moves to B8 (no longer the
# tab shown), Pet's selection is left untouched.

$wb = $excel.ActiveWorkbook

$wsGrocery = $wb.Worksheets.Item("Grocery")
$wsPet     = $wb.Worksheets.Item("Pet")
$wsBath    = $wb.Worksheets.Item("Bath")

# --- Grocery: Fruits, Vegetables, Basic Cat Food, Lotion, Tea -> 3 ---
$wsGrocery.Range("B3").Value = 3
$wsGrocery.Range("B4").Value = 3
$wsGrocery.Range("B5").Value = 3
$wsGrocery.Range("B7").Value = 3
$wsGrocery.Range("B8").Value = 3

# --- Pet: Premium Cat Food, Premium Dog Food -> 3 (Catnip untouched) ---
$wsPet.Range("B3").Value = 3
$wsPet.Range("B4").Value = 3

# --- Bath: Face Scrub, Shampoo Bar, Hair Conditioner Balm -> 3 ---
$wsBath.Range("B4").Value = 3
$wsBath.Range("B5").Value = 3
$wsBath.Range("B6").Value = 3

# --- Selection / active-sheet bookkeeping ---
$wsGrocery.Range("B8").Select()

$wsBath.Activate()
$wsBath.Range("B6").Select()
